$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: iaest-measure:* -> iaest-dimension:* / sdmx-dimension:refArea
$ws.Range("A3").Value = "iaest-dimension:edad-grandes-grupos"
$ws.Range("E3").Value = "sdmx-dimension:refArea"
$ws.Range("F3").Value = "sdmx-dimension:refArea"
$ws.Range("G3").Value = "sdmx-dimension:refArea"
$ws.Range("I3").Value = "iaest-dimension:sexo"

# Row 4: medida -> dim for the columns that became dimensions (H4 was already "dim")
$ws.Range("A4").Value = "dim"
$ws.Range("E4").Value = "dim"
$ws.Range("F4").Value = "dim"
$ws.Range("G4").Value = "dim"
$ws.Range("I4").Value = "dim"

# Row 5: xsd:string -> skos:Concept / URI-Comunidad / URI-Provincia
$ws.Range("A5").Value = "skos:Concept"
$ws.Range("E5").Value = "URI-Comunidad"
$ws.Range("F5").Value = "URI-Provincia"
$ws.Range("G5").Value = "URI-Provincia"
$ws.Range("I5").Value = "skos:Concept"

# New row 6: mapping file references
# copy the existing row's style onto the brand-new cells before writing values,
# so they pick up the same cell style (s="1") as the rest of the sheet.
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("I5").Copy()
$ws.Range("I6").PasteSpecial(-4122)

$ws.Range("A6").Value = "mapping-edad-grandes-grupos.xlsx"
$ws.Range("I6").Value = "mapping-sexo.xlsx"
